$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row text ---------------------------------------------
# Title (merged A1:K1 -> A1:J1 after column removal)
$ws.Range("A1").Value = "Matriz Integrada de Compromisos de la Unidad de Multiunidad"

# Column B no longer "Origen compromiso" -> becomes "Codigo global"
$ws.Range("B2").Value = "Código global"
# Column C no longer "Codigo global" -> becomes "Capitulo"
$ws.Range("C2").Value = "Capítulo"
# D (Seccion) and E (Pagina) stay the same
# Column F no longer "Aspecto asociado" -> becomes "Temporalidad"
$ws.Range("F2").Value = "Temporalidad"
# Column G no longer "Contenido original del compromiso" -> becomes "Fecha de inicio"
$ws.Range("G2").Value = "Fecha de inicio"
# Column H no longer "Nombre de revisor" -> becomes "Frecuencia"
$ws.Range("H2").Value = "Frecuencia"
# Column I no longer "Notas adicionales" -> becomes "Criticidad"
$ws.Range("I2").Value = "Criticidad"
# Column J no longer "Operacion" -> becomes "Estado de cumplimiento"
$ws.Range("J2").Value = "Estado de cumplimiento"

# --- Remove the trailing column K (was "Cierre") -------------------------
$ws.Range("K:K").Delete()

# --- Fix up the workbook-level defined names ------------------------------
# nrocorrelativo (A1:A2) stays as-is.
$wb.Names.Item("origencompromiso").Delete()
$wb.Names.Item("codigoglobal").Delete()
$wb.Names.Item("aspambasoc").Delete()
$wb.Names.Item("contorigcomp").Delete()
$wb.Names.Item("nombrerevisor").Delete()
$wb.Names.Item("comentarios").Delete()
$wb.Names.Item("operacion").Delete()
$wb.Names.Item("cierre").Delete()
# seccion, pagina stay as-is.

$wb.Names.Add("codigoglobal", "=Compromisos!`$B`$1:`$B`$2")
$wb.Names.Add("capitulo", "=Compromisos!`$C`$1:`$C`$2")
$wb.Names.Add("temporalidad", "=Compromisos!`$F`$1:`$F`$2")
$wb.Names.Add("fechainicio", "=Compromisos!`$G`$1:`$G`$2")
$wb.Names.Add("frecuencia", "=Compromisos!`$H`$1:`$H`$2")
$wb.Names.Add("criticidad", "=Compromisos!`$I`$1:`$I`$2")
$wb.Names.Add("estadocumplimiento", "=Compromisos!`$J`$1:`$J`$2")
